# Scheduled-runner update: refresh computed profit/price columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the latest
# market-board derived figures.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 10458.385
$ws.Cells.Item(34, 9).Value = 10458.385
$ws.Cells.Item(34, 11).Value = 10458.385
$ws.Cells.Item(34, 13).Value = -10255.385
$ws.Cells.Item(36, 8).Value = 10458.385
$ws.Cells.Item(36, 9).Value = 10458.385
$ws.Cells.Item(36, 11).Value = 10458.385
$ws.Cells.Item(36, 13).Value = -9743.385
$ws.Cells.Item(82, 8).Value = 5125.3335
$ws.Cells.Item(82, 9).Value = 5125.3335
$ws.Cells.Item(82, 11).Value = 15376.0005
$ws.Cells.Item(82, 13).Value = -14970.0005
$ws.Cells.Item(85, 8).Value = 5125.3335
$ws.Cells.Item(85, 9).Value = 5125.3335
$ws.Cells.Item(85, 11).Value = 15376.0005
$ws.Cells.Item(85, 13).Value = -13972.0005
$ws.Cells.Item(100, 8).Value = 2321.2144
$ws.Cells.Item(100, 9).Value = 1995.5454
$ws.Cells.Item(100, 10).Value = 3515.3333
$ws.Cells.Item(100, 11).Value = 1995.5454
$ws.Cells.Item(100, 12).Value = 3515.3333
$ws.Cells.Item(100, 13).Value = -1454.5454
$ws.Cells.Item(100, 14).Value = -4597.3333
$ws.Cells.Item(103, 8).Value = 248.5
$ws.Cells.Item(103, 9).Value = 248.5
$ws.Cells.Item(103, 11).Value = 745.5
$ws.Cells.Item(103, 13).Value = -159.5
$ws.Cells.Item(106, 8).Value = 1140.8572
$ws.Cells.Item(106, 9).Value = 983.4666999999999
$ws.Cells.Item(106, 11).Value = 983.4666999999999
$ws.Cells.Item(106, 13).Value = -352.4666999999999
$ws.Cells.Item(108, 8).Value = 84992.5
$ws.Cells.Item(108, 10).Value = 84992.5
$ws.Cells.Item(108, 12).Value = 84992.5
$ws.Cells.Item(108, 14).Value = -92672.5
$ws.Cells.Item(109, 8).Value = 64455
$ws.Cells.Item(109, 10).Value = 64455
$ws.Cells.Item(109, 12).Value = 64455
$ws.Cells.Item(109, 14).Value = -67229
$ws.Cells.Item(110, 8).Value = 60897.6
$ws.Cells.Item(110, 10).Value = 60897.6
$ws.Cells.Item(110, 12).Value = 60897.6
$ws.Cells.Item(110, 14).Value = -69077.60000000001
$ws.Cells.Item(112, 8).Value = 1312.76
$ws.Cells.Item(112, 10).Value = 1314.5
$ws.Cells.Item(112, 12).Value = 3943.5
$ws.Cells.Item(112, 14).Value = -6159.5
$ws.Cells.Item(117, 8).Value = 81476.5
$ws.Cells.Item(117, 10).Value = 81476.5
$ws.Cells.Item(117, 12).Value = 81476.5
$ws.Cells.Item(117, 14).Value = -90654.5
$ws.Cells.Item(120, 8).Value = 49495
$ws.Cells.Item(120, 10).Value = 49495
$ws.Cells.Item(120, 12).Value = 49495
$ws.Cells.Item(120, 14).Value = -59171
$ws.Cells.Item(123, 8).Value = 64620.375
$ws.Cells.Item(123, 10).Value = 64620.375
$ws.Cells.Item(123, 12).Value = 64620.375
$ws.Cells.Item(123, 14).Value = -74420.375
$ws.Cells.Item(132, 8).Value = 1509
$ws.Cells.Item(132, 9).Value = 1387.9387
$ws.Cells.Item(132, 11).Value = 4163.8161
$ws.Cells.Item(132, 13).Value = -1633.8161
$ws.Cells.Item(138, 8).Value = 2251.6428
$ws.Cells.Item(138, 9).Value = 1534.8334
$ws.Cells.Item(138, 10).Value = 2789.25
$ws.Cells.Item(138, 11).Value = 4604.5002
$ws.Cells.Item(138, 12).Value = 8367.75
$ws.Cells.Item(138, 13).Value = 535.4997999999996
$ws.Cells.Item(138, 14).Value = -18647.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 27482.625
$ws.Cells.Item(7, 10).Value = 27482.625
$ws.Cells.Item(7, 12).Value = 27482.625
$ws.Cells.Item(7, 14).Value = -27710.625
$ws.Cells.Item(32, 8).Value = 3642.923
$ws.Cells.Item(32, 9).Value = 3022.9333
$ws.Cells.Item(32, 11).Value = 3022.9333
$ws.Cells.Item(32, 13).Value = -2735.9333
$ws.Cells.Item(35, 8).Value = 1634.5
$ws.Cells.Item(35, 9).Value = 1634.5
$ws.Cells.Item(35, 11).Value = 1634.5
$ws.Cells.Item(35, 13).Value = -1228.5
$ws.Cells.Item(74, 8).Value = 33631.773
$ws.Cells.Item(74, 9).Value = 36761.965
$ws.Cells.Item(74, 11).Value = 36761.965
$ws.Cells.Item(74, 13).Value = -35887.965
$ws.Cells.Item(77, 8).Value = 33631.773
$ws.Cells.Item(77, 9).Value = 36761.965
$ws.Cells.Item(77, 11).Value = 183809.825
$ws.Cells.Item(77, 13).Value = -179441.825
$ws.Cells.Item(104, 8).Value = 29859.5
$ws.Cells.Item(104, 10).Value = 29859.5
$ws.Cells.Item(104, 12).Value = 29859.5
$ws.Cells.Item(104, 14).Value = -36847.5
$ws.Cells.Item(118, 8).Value = 57017.7
$ws.Cells.Item(118, 10).Value = 57017.7
$ws.Cells.Item(118, 12).Value = 57017.7
$ws.Cells.Item(118, 14).Value = -60331.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(13, 8).Value = 86988
$ws.Cells.Item(13, 10).Value = 86988
$ws.Cells.Item(13, 12).Value = 86988
$ws.Cells.Item(13, 14).Value = -87324
$ws.Cells.Item(53, 8).Value = 24998
$ws.Cells.Item(53, 10).Value = 24998
$ws.Cells.Item(53, 12).Value = 24998
$ws.Cells.Item(53, 14).Value = -26146
$ws.Cells.Item(108, 8).Value = 89994
$ws.Cells.Item(108, 10).Value = 89994
$ws.Cells.Item(108, 12).Value = 89994
$ws.Cells.Item(108, 14).Value = -97674
$ws.Cells.Item(109, 8).Value = 99989.2
$ws.Cells.Item(109, 10).Value = 99989.2
$ws.Cells.Item(109, 12).Value = 99989.2
$ws.Cells.Item(109, 14).Value = -102763.2
$ws.Cells.Item(114, 8).Value = 89989.2
$ws.Cells.Item(114, 10).Value = 89989.2
$ws.Cells.Item(114, 12).Value = 89989.2
$ws.Cells.Item(114, 14).Value = -98667.2
$ws.Cells.Item(118, 8).Value = 54535.5
$ws.Cells.Item(118, 10).Value = 54198.668
$ws.Cells.Item(118, 12).Value = 54198.668
$ws.Cells.Item(118, 14).Value = -57512.668
$ws.Cells.Item(122, 8).Value = 95439.5
$ws.Cells.Item(122, 10).Value = 95439.5
$ws.Cells.Item(122, 12).Value = 95439.5
$ws.Cells.Item(122, 14).Value = -105239.5
$ws.Cells.Item(132, 8).Value = 27797.684
$ws.Cells.Item(132, 10).Value = 27947.299
$ws.Cells.Item(132, 12).Value = 27947.299
$ws.Cells.Item(132, 14).Value = -38067.299
$ws.Cells.Item(135, 8).Value = 95267.5
$ws.Cells.Item(135, 10).Value = 95267.5
$ws.Cells.Item(135, 12).Value = 95267.5
$ws.Cells.Item(135, 14).Value = -105407.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2827.5122
$ws.Cells.Item(31, 9).Value = 2054.389
$ws.Cells.Item(31, 10).Value = 3432.5652
$ws.Cells.Item(31, 11).Value = 2054.389
$ws.Cells.Item(31, 12).Value = 3432.5652
$ws.Cells.Item(31, 13).Value = -1759.389
$ws.Cells.Item(31, 14).Value = -4022.5652
$ws.Cells.Item(34, 8).Value = 2827.5122
$ws.Cells.Item(34, 9).Value = 2054.389
$ws.Cells.Item(34, 10).Value = 3432.5652
$ws.Cells.Item(34, 11).Value = 2054.389
$ws.Cells.Item(34, 12).Value = 3432.5652
$ws.Cells.Item(34, 13).Value = -1852.389
$ws.Cells.Item(34, 14).Value = -3836.5652
$ws.Cells.Item(58, 8).Value = 1334.3125
$ws.Cells.Item(58, 9).Value = 1097.8572
$ws.Cells.Item(58, 10).Value = 1785.7273
$ws.Cells.Item(58, 11).Value = 1097.8572
$ws.Cells.Item(58, 12).Value = 1785.7273
$ws.Cells.Item(58, 13).Value = -894.8571999999999
$ws.Cells.Item(58, 14).Value = -2191.7273
$ws.Cells.Item(114, 8).Value = 53842.625
$ws.Cells.Item(114, 10).Value = 53842.625
$ws.Cells.Item(114, 12).Value = 53842.625
$ws.Cells.Item(114, 14).Value = -62520.625
$ws.Cells.Item(116, 8).Value = 49119.25
$ws.Cells.Item(116, 10).Value = 49119.25
$ws.Cells.Item(116, 12).Value = 49119.25
$ws.Cells.Item(116, 14).Value = -58297.25
$ws.Cells.Item(118, 8).Value = 58096.555
$ws.Cells.Item(118, 10).Value = 58096.555
$ws.Cells.Item(118, 12).Value = 58096.555
$ws.Cells.Item(118, 14).Value = -61410.555
$ws.Cells.Item(119, 8).Value = 65203
$ws.Cells.Item(119, 10).Value = 65203
$ws.Cells.Item(119, 12).Value = 65203
$ws.Cells.Item(119, 14).Value = -74879
$ws.Cells.Item(136, 8).Value = 1334.3125
$ws.Cells.Item(136, 9).Value = 1097.8572
$ws.Cells.Item(136, 10).Value = 1785.7273
$ws.Cells.Item(136, 11).Value = 3293.5716
$ws.Cells.Item(136, 12).Value = 5357.1819
$ws.Cells.Item(136, 13).Value = -743.5715999999998
$ws.Cells.Item(136, 14).Value = -10457.1819
$ws.Cells.Item(138, 8).Value = 89528.086
$ws.Cells.Item(138, 10).Value = 91693.45
$ws.Cells.Item(138, 12).Value = 91693.45
$ws.Cells.Item(138, 14).Value = -101973.45

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 23.333334
$ws.Cells.Item(19, 9).Value = 22.5
$ws.Cells.Item(19, 10).Value = 25
$ws.Cells.Item(19, 11).Value = 67.5
$ws.Cells.Item(19, 12).Value = 75
$ws.Cells.Item(19, 13).Value = 106.5
$ws.Cells.Item(19, 14).Value = -423
$ws.Cells.Item(45, 8).Value = 4165
$ws.Cells.Item(45, 10).Value = 4165
$ws.Cells.Item(45, 12).Value = 12495
$ws.Cells.Item(45, 14).Value = -13559
$ws.Cells.Item(113, 8).Value = 112329.555

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 637836.8
$ws.Cells.Item(3, 9).Value = 200
$ws.Cells.Item(3, 10).Value = 701600.5
$ws.Cells.Item(3, 11).Value = 200
$ws.Cells.Item(3, 12).Value = 701600.5
$ws.Cells.Item(3, 13).Value = -84
$ws.Cells.Item(3, 14).Value = -701832.5
$ws.Cells.Item(114, 8).Value = 71495.71000000001
$ws.Cells.Item(114, 10).Value = 71495.71000000001
$ws.Cells.Item(114, 12).Value = 71495.71000000001
$ws.Cells.Item(114, 14).Value = -80173.71000000001
$ws.Cells.Item(140, 8).Value = 94184.42999999999
$ws.Cells.Item(140, 10).Value = 94184.42999999999
$ws.Cells.Item(140, 12).Value = 94184.42999999999
$ws.Cells.Item(140, 14).Value = -104544.43

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 5499
$ws.Cells.Item(11, 9).Value = 5499
$ws.Cells.Item(11, 11).Value = 5499
$ws.Cells.Item(11, 13).Value = -5359
$ws.Cells.Item(132, 8).Value = 1754.4546
$ws.Cells.Item(132, 9).Value = 1427.65
$ws.Cells.Item(132, 11).Value = 4282.950000000001
$ws.Cells.Item(132, 13).Value = -1752.950000000001
$ws.Cells.Item(136, 8).Value = 2727.1892
$ws.Cells.Item(136, 9).Value = 2623.5862
$ws.Cells.Item(136, 11).Value = 7870.758600000001
$ws.Cells.Item(136, 13).Value = -5320.758600000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2175111.5
$ws.Cells.Item(132, 9).Value = 988.4375
$ws.Cells.Item(132, 11).Value = 2965.3125
$ws.Cells.Item(132, 13).Value = -435.3125
$ws.Cells.Item(136, 8).Value = 1992
$ws.Cells.Item(136, 9).Value = 1805
$ws.Cells.Item(136, 11).Value = 5415
$ws.Cells.Item(136, 13).Value = -2865
